# Regenerate save_data: replace column G ("K") values (previously derived
# from a "Strike#" style computation) with the newly recomputed s_vals.
# Only the K column (G) values for the data rows (2-33) change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values, keyed by row number (row 2 = first data row / index 0 ... row 33 = index 31)
$kValues = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 1
    6  = 1
    7  = 1
    8  = 0
    9  = 0
    10 = 0
    11 = 1
    12 = 2
    13 = 3
    14 = 1
    15 = 1
    16 = 2
    17 = 4
    18 = 1
    19 = 0
    20 = 0
    21 = 3
    22 = 0
    23 = 1
    24 = 2
    25 = 1
    26 = 0
    27 = 0
    28 = 0
    29 = 0
    30 = 1
    31 = 4
    32 = 1
    33 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
